$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts rows 8+ down by one),
# reflecting the new cbt02_2 / cbt02_3 split and the appended cbt02_6 entry.
$ws.Rows(8).Insert()

# Rewrite row content (cols A and B) for rows 3 through 82 to match the revised
# question ordering / timing.
$ws.Cells.Item(3, 1).Value = "For example, if you feel depressed, you will likely have many negative thoughts, possibly about yourself or others, or about your life in general."
$ws.Cells.Item(3, 2).Value = "cbt01_2"
$ws.Cells.Item(4, 1).Value = "Also, when you're depressed, you're more likely to act in ways that feed your negative thoughts and sadness."
$ws.Cells.Item(4, 2).Value = "cbt01_3"
$ws.Cells.Item(5, 1).Value = "These behaviors may include withdrawing socially, avoiding tasks, and poor self-care."
$ws.Cells.Item(5, 2).Value = "cbt01_4"
$ws.Cells.Item(6, 1).Value = "In therapy, you will learn skills to change unhelpful thought and behavior patterns, which will help you feel better emotionally."
$ws.Cells.Item(6, 2).Value = "cbt01_5"
$ws.Cells.Item(7, 1).Value = "When depressed, people frequently overlook their negative thoughts or cognitions, and assume what they are thinking is truth or reality. In therapy, we work to question the reality created by depressive thoughts."
$ws.Cells.Item(7, 2).Value = "cbt02_1"
$ws.Cells.Item(8, 1).Value = "The first step in doing this is learning to catch or identify automatic thoughts."
$ws.Cells.Item(8, 2).Value = "cbt02_2"
$ws.Cells.Item(9, 1).Value = "Automatic thoughts are thoughts that pop into our heads, and most of the time, they happen automatically and outside of our awareness."
$ws.Cells.Item(9, 2).Value = "cbt02_3"
$ws.Cells.Item(10, 1).Value = "We are usually more aware of the emotions these thoughts make us feel, such as sad or depressed, and we do not question what we are thinking. But often, the thoughts we are having are distorting reality in some way."
$ws.Cells.Item(10, 2).Value = "cbt02_4"
$ws.Cells.Item(11, 1).Value = "Before we can even start questioning our thoughts, we have to slow down our thoughts and know what may be fueling our sadness."
$ws.Cells.Item(11, 2).Value = "cbt02_5"
$ws.Cells.Item(12, 1).Value = "One way to do this is to systematically write down moments we feel a negative emotion, and identify what thoughts we were having during those moments."
$ws.Cells.Item(12, 2).Value = "cbt02_6"
$ws.Cells.Item(13, 1).Value = "As you get better at identifying your automatic thoughts, you will see that these negative thoughts that flood your mind can fuel your sadness and depression."
$ws.Cells.Item(13, 2).Value = "cbt03_1"
$ws.Cells.Item(14, 1).Value = "These thoughts can keep you feeling lethargic and inadequate."
$ws.Cells.Item(14, 2).Value = "cbt03_2"
$ws.Cells.Item(15, 1).Value = "Finding ways to notice and evaluate the accuracy of your negative thoughts are key to feeling better."
$ws.Cells.Item(15, 2).Value = "cbt03_3"
$ws.Cells.Item(16, 1).Value = "When we look at your thoughts, sometimes they are true, not true, or have a grain of truth."
$ws.Cells.Item(16, 2).Value = "cbt04_1"
$ws.Cells.Item(17, 1).Value = "We will teach you how to evaluate the accuracy of your thoughts by looking at the evidence."
$ws.Cells.Item(17, 2).Value = "cbt04_2"
$ws.Cells.Item(18, 1).Value = "What evidence is there that the thought is true?"
$ws.Cells.Item(18, 2).Value = "cbt04_3"
$ws.Cells.Item(19, 1).Value = "What evidence is there that the thought is not true?"
$ws.Cells.Item(19, 2).Value = "cbt04_4"
$ws.Cells.Item(20, 1).Value = "Because your automatic thoughts are so quick, they are hard to notice or remember."
$ws.Cells.Item(20, 2).Value = "cbt05_1"
$ws.Cells.Item(21, 1).Value = "We want to practice identifying them by writing them down on what we call a ""Thought Record."""
$ws.Cells.Item(21, 2).Value = "cbt05_2"
$ws.Cells.Item(22, 1).Value = "We usually notice how we're feeling before we notice what we're thinking."
$ws.Cells.Item(22, 2).Value = "cbt05_3"
$ws.Cells.Item(23, 1).Value = "I want you to notice moments when you feel a negative emotion and write down on the Thought Record what the situation was, what emotions you were feeling, and then what thoughts you had when the event or situation happened."
$ws.Cells.Item(23, 2).Value = "cbt05_4"
$ws.Cells.Item(24, 1).Value = "It is crucial to write down your automatic thoughts."
$ws.Cells.Item(24, 2).Value = "cbt06_1"
$ws.Cells.Item(25, 1).Value = "Writing them down forces you to look at them more objectively."
$ws.Cells.Item(25, 2).Value = "cbt06_2"
$ws.Cells.Item(26, 1).Value = "It also helps you see if there are any patterns or mistakes in your thoughts."
$ws.Cells.Item(26, 2).Value = "cbt06_3"
$ws.Cells.Item(27, 1).Value = "There are many common mistakes, or “cognitive distortions,” people make in their thinking."
$ws.Cells.Item(27, 2).Value = "cbt07_1"
$ws.Cells.Item(28, 1).Value = "By figuring out what cognitive distortions you may be having, you can respond to situations in a more healthy way. "
$ws.Cells.Item(28, 2).Value = "cbt07_2"
$ws.Cells.Item(29, 1).Value = "One cognitive distortion is called ""all or nothing thinking,"" where you see things in very black and white terms."
$ws.Cells.Item(29, 2).Value = "cbt07_3"
$ws.Cells.Item(30, 1).Value = "For example, if you were feeling depressed and unable to finish your homework, you might think that you are a failure and will be a failure forever."
$ws.Cells.Item(30, 2).Value = "cbt07_4"
$ws.Cells.Item(31, 1).Value = "Earlier we talked about finding evidence for your automatic thoughts."
$ws.Cells.Item(31, 2).Value = "cbt08_1"
$ws.Cells.Item(32, 1).Value = "Doing a ""behavioral experiment"" is one way to collect and examine the evidence. To do this, you may act out your automatic thought to see if what you predict happens. "
$ws.Cells.Item(32, 2).Value = "cbt08_2"
$ws.Cells.Item(33, 1).Value = "For example, say you think you will feel worse if you get out of bed and take a shower because you feel so sad."
$ws.Cells.Item(33, 2).Value = "cbt08_3"
$ws.Cells.Item(34, 1).Value = "You can test that out by getting out of bed and showering when you feel depressed and want to lay in bed. "
$ws.Cells.Item(34, 2).Value = "cbt08_4"
$ws.Cells.Item(35, 1).Value = "We can then see if you do indeed feel worse as your automatic thought predicts you will."
$ws.Cells.Item(35, 2).Value = "cbt08_5"
$ws.Cells.Item(36, 1).Value = "When we do behavioral experiments, we will create hypotheses to test. In “""ypothesis testing,"" we challenge unhelpful behaviors and thoughts by testing out alternative ways of responding. "
$ws.Cells.Item(36, 2).Value = "cbt09_1"
$ws.Cells.Item(37, 1).Value = "We do this by developing a behavioral experiment to test if your new way of thinking or new behavior makes you feel better than your current thoughts and behaviors. We then ""collect data"" as you run the experiment."
$ws.Cells.Item(37, 2).Value = "cbt09_2"
$ws.Cells.Item(38, 1).Value = "For example, you know that avoiding your email is unhelpful for your depression. You have the hypothesis that checking your email will only worsen your mood and depression."
$ws.Cells.Item(38, 2).Value = "cbt09_3"
$ws.Cells.Item(39, 1).Value = "You can then test out this hypothesis by checking your email at a specific time, and rating your mood before and after you checked your email."
$ws.Cells.Item(39, 2).Value = "cbt09_4"
$ws.Cells.Item(40, 1).Value = "Once you've collected the data, we will evaluate the benefits of challenging your avoidance."
$ws.Cells.Item(40, 2).Value = "cbt09_5"
$ws.Cells.Item(41, 1).Value = "Doing these behavioral experiments is difficult and can bring up many emotions."
$ws.Cells.Item(41, 2).Value = "cbt10_1"
$ws.Cells.Item(42, 1).Value = "To help better understand how intensely you are feeling these emotions, we will teach you to use a rating scale called the subjective units of distress scale (SUDS)."
$ws.Cells.Item(42, 2).Value = "cbt10_2"
$ws.Cells.Item(43, 1).Value = "Using the SUDS, your emotions will range in intensity from 0 (not noticeable) to 100 (the highest extreme)."
$ws.Cells.Item(43, 2).Value = "cbt10_3"
$ws.Cells.Item(44, 1).Value = "The SUDS will help you communicate what you are feeling more accurately and in an understandable way to both of us."
$ws.Cells.Item(44, 2).Value = "cbt10_4"
$ws.Cells.Item(45, 1).Value = "As you become more skilled at using the SUDS, you will feel clearer about your feelings and more able to make decisions about what you would like to try to do in therapy."
$ws.Cells.Item(45, 2).Value = "cbt10_5"
$ws.Cells.Item(46, 1).Value = "Effectively managing stressful life problems requires a planful approach."
$ws.Cells.Item(46, 2).Value = "pst01_1"
$ws.Cells.Item(47, 1).Value = "The set of skills required to do this require both learning and practice and include four planful problem-solving skills."
$ws.Cells.Item(47, 2).Value = "pst01_2"
$ws.Cells.Item(48, 1).Value = "We teach you these skills and help you practice them in your everyday life by completing PST planful problem-solving worksheets."
$ws.Cells.Item(48, 2).Value = "pst01_3"
$ws.Cells.Item(49, 1).Value = "The four skills are problem definition, generating alternatives, decision-making, and solution implementation and verification."
$ws.Cells.Item(49, 2).Value = "pst01_4"
$ws.Cells.Item(50, 1).Value = "One skill that may help you is called ""problem definition,"" or the clarifying the nature of a problem."
$ws.Cells.Item(50, 2).Value = "pst02_1"
$ws.Cells.Item(51, 1).Value = "You will learn how to set a realistic problem-solving goal and identify the obstacles that are currently preventing you from reaching that goal."
$ws.Cells.Item(51, 2).Value = "pst02_2"
$ws.Cells.Item(52, 1).Value = "For example, if you want to raise your GPA, we will figure out a realistic GPA that you can achieve and identify the very real obstacles that are currently preventing you from reaching that."
$ws.Cells.Item(52, 2).Value = "pst02_3"
$ws.Cells.Item(53, 1).Value = "One skill that may help you is the ""generation of alternatives,"" in which we teach you to use your creative skills to brainstorm different types of solutions."
$ws.Cells.Item(53, 2).Value = "pst03_1"
$ws.Cells.Item(54, 1).Value = "You will come up with as many solutions as you can. Not all will be the best solution, but it helps you think of alternatives."
$ws.Cells.Item(54, 2).Value = "pst03_2"
$ws.Cells.Item(55, 1).Value = "For example, if you are struggling to talk to your roommate about cleaning up after themselves, you can use the brainstorming tool to discover that there are many ways to get closer to your goal and get past the obstacles in your way."
$ws.Cells.Item(55, 2).Value = "pst03_3"
$ws.Cells.Item(56, 1).Value = "Using the skill of ""decision making,"" you will learn how to look at the likely consequences of different solution ideas."
$ws.Cells.Item(56, 2).Value = "pst04_1"
$ws.Cells.Item(57, 1).Value = "After looking at these consequences, you will learn how to develop an action plan that is geared toward achieving the problem-solving goal."
$ws.Cells.Item(57, 2).Value = "pst04_2"
$ws.Cells.Item(58, 1).Value = "For example, if you want to raise your GPA, you can develop an action plan that represents the best solution match for you."
$ws.Cells.Item(58, 2).Value = "pst04_3"
$ws.Cells.Item(59, 1).Value = "One skill that may help you is called ""solution implementation and verification."""
$ws.Cells.Item(59, 2).Value = "pst05_1"
$ws.Cells.Item(60, 1).Value = "This skill involves carrying out the action plan, monitoring and evaluating the consequences of the plan, and determining whether one's problem-solving efforts have been successful."
$ws.Cells.Item(60, 2).Value = "pst05_2"
$ws.Cells.Item(61, 1).Value = "For example, if you are struggling to talk to your roommate about cleaning up after themselves, once you have figured out the best solution, you can carry it out, see what happens, and evaluate if the outcome is what you wanted."
$ws.Cells.Item(61, 2).Value = "pst05_3"
$ws.Cells.Item(62, 1).Value = "Through ""externalization"" you will experience how writing things down, recording messages for yourself on your iPhone, or talking through a difficult problem helps you to be less overwhelmed."
$ws.Cells.Item(62, 2).Value = "pst06_1"
$ws.Cells.Item(63, 1).Value = "Try getting the brain overload you are experiencing with a difficult problem ""out of your head"" and onto paper, and you may notice that the thinking part of your brain can better understand and begin to organize this challenging problem or goal."
$ws.Cells.Item(63, 2).Value = "pst06_2"
$ws.Cells.Item(64, 1).Value = "For example, trying to talk to your parents about your current relationship may be really hard."
$ws.Cells.Item(64, 2).Value = "pst06_3"
$ws.Cells.Item(65, 1).Value = "Let's write down all of the thoughts, feelings, and concerns that are contributing to your experience of being overwhelmed as a way to begin to organize all of this information."
$ws.Cells.Item(65, 2).Value = "pst06_4"
$ws.Cells.Item(66, 1).Value = "Through visualization, you learn to use visual imagery to help understand and clarify a current problem or goal."
$ws.Cells.Item(66, 2).Value = "pst07_1"
$ws.Cells.Item(67, 1).Value = "One way is to try picturing the problem in your imagination to help you better define it."
$ws.Cells.Item(67, 2).Value = "pst07_2"
$ws.Cells.Item(68, 1).Value = "You can also use visualization to rehearse how you'll carry out a solution or action plan that you have developed."
$ws.Cells.Item(68, 2).Value = "pst07_3"
$ws.Cells.Item(69, 1).Value = "Finally, you can use visualization to help calm you when you experience strong emotions associated with stress."
$ws.Cells.Item(69, 2).Value = "pst07_4"
$ws.Cells.Item(70, 1).Value = "We teach you the use of ""simplification"" to break down a large or complex problem to make it more manageable."
$ws.Cells.Item(70, 2).Value = "pst08_1"
$ws.Cells.Item(71, 1).Value = "You will learn how to break down these big problems into smaller steps."
$ws.Cells.Item(71, 2).Value = "pst08_2"
$ws.Cells.Item(72, 1).Value = "For example, consider how you would begin to break down a complex situation like getting into medical school into smaller pieces to accomplish one at a time."
$ws.Cells.Item(72, 2).Value = "pst08_3"
$ws.Cells.Item(73, 1).Value = "When facing a stressful problem or daunting goal, it can be hard to believe that there are ways to effectively manage the stress, reach a goal, or solve the problem."
$ws.Cells.Item(73, 2).Value = "pst09_1"
$ws.Cells.Item(74, 1).Value = "One tool that can help uses visualization in a special way to give you a glimpse of the future."
$ws.Cells.Item(74, 2).Value = "pst09_2"
$ws.Cells.Item(75, 1).Value = "Try imagining what it would be like at a moment in time in the future when a stressful problem you are facing is largely resolved and the obstacles overcome, such as making a change in your plan for a career."
$ws.Cells.Item(75, 2).Value = "pst09_3"
$ws.Cells.Item(76, 1).Value = "This won't immediately solve the problem but will give you an experience of what it would feel like to reach your goal and experience a ""light at the end of the tunnel."""
$ws.Cells.Item(76, 2).Value = "pst09_4"
$ws.Cells.Item(77, 1).Value = "We all need a picture in our head of what we are working toward. We will continue to practice doing this together."
$ws.Cells.Item(77, 2).Value = "pst09_5"
$ws.Cells.Item(78, 1).Value = "This important toolkit is referred to as the ""SSTA"" method."
$ws.Cells.Item(78, 2).Value = "pst10_1"
$ws.Cells.Item(79, 1).Value = "In SSTA, the first S = Stop (and be aware of what you are experiencing)."
$ws.Cells.Item(79, 2).Value = "pst10_2"
$ws.Cells.Item(80, 1).Value = "The second S = Slow down and take a moment to ""turn down the volume of strong emotions."" You can take a few deep breaths, or use other techniques, so you can still listen to your feelings which give you important information, but allow your brain to keep working."
$ws.Cells.Item(80, 2).Value = "pst10_3"
$ws.Cells.Item(81, 1).Value = "Ultimately, you will need to think carefully and planfully about an action plan that gives you the best chance of reaching your goals or solving a problem. Only after reducing intense emotional arousal and ""turning down the volume,"" can you planfully and carefully."
$ws.Cells.Item(81, 2).Value = "pst10_4"
$ws.Cells.Item(82, 1).Value = "Finally, T = Think, and A = Act."
$ws.Cells.Item(82, 2).Value = "pst10_5"

# Update view to reflect the new scroll position / selection used when the
# author finished this edit.
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("A58").Select()
